$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 6,7,8: cyclic reassignment (6<-7old, 7<-8old, 8<-6old) ---
$ws.Cells.Item(6, 6).Value = 'Wisla Pulawy'
$ws.Cells.Item(6, 7).Value = 1
$ws.Cells.Item(6, 8).Value = 'Chojniczanka'
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 10).Value = 2.18
$ws.Cells.Item(6, 11).Value = '22/07/2023 14:42'
$ws.Cells.Item(6, 12).Value = 2.46
$ws.Cells.Item(6, 13).Value = '22/07/2023 17:06'
$ws.Cells.Item(6, 14).Value = 3.31
$ws.Cells.Item(6, 15).Value = '22/07/2023 14:42'
$ws.Cells.Item(6, 16).Value = 3.45
$ws.Cells.Item(6, 17).Value = '22/07/2023 17:35'
$ws.Cells.Item(6, 18).Value = 3.05
$ws.Cells.Item(6, 19).Value = '22/07/2023 14:42'
$ws.Cells.Item(6, 20).Value = 2.6
$ws.Cells.Item(6, 21).Value = '22/07/2023 17:06'
$ws.Cells.Item(6, 22).Value = 'https://www.betexplorer.com/football/poland/division-2/wisla-pulawy-chojniczanka/rc7er10t/'
$ws.Cells.Item(7, 6).Value = 'Ol. Grudziadz'
$ws.Cells.Item(7, 7).Value = 3
$ws.Cells.Item(7, 8).Value = 'Sandecja Nowy S.'
$ws.Cells.Item(7, 9).Value = 1
$ws.Cells.Item(7, 10).Value = 2.49
$ws.Cells.Item(7, 11).Value = '22/07/2023 14:42'
$ws.Cells.Item(7, 12).Value = 2.39
$ws.Cells.Item(7, 13).Value = '22/07/2023 17:36'
$ws.Cells.Item(7, 14).Value = 3.26
$ws.Cells.Item(7, 15).Value = '22/07/2023 14:42'
$ws.Cells.Item(7, 16).Value = 3.35
$ws.Cells.Item(7, 17).Value = '22/07/2023 17:36'
$ws.Cells.Item(7, 18).Value = 2.64
$ws.Cells.Item(7, 19).Value = '22/07/2023 14:42'
$ws.Cells.Item(7, 20).Value = 2.82
$ws.Cells.Item(7, 21).Value = '22/07/2023 17:36'
$ws.Cells.Item(7, 22).Value = 'https://www.betexplorer.com/football/poland/division-2/ol-grudziadz-sandecja-nowy-s/GUiPyHqO/'
$ws.Cells.Item(8, 6).Value = 'GKS Jastrzebie'
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = 'S. Wola'
$ws.Cells.Item(8, 9).Value = 0
$ws.Cells.Item(8, 10).Value = 2.57
$ws.Cells.Item(8, 11).Value = '22/07/2023 14:42'
$ws.Cells.Item(8, 12).Value = 2.64
$ws.Cells.Item(8, 13).Value = '22/07/2023 17:59'
$ws.Cells.Item(8, 14).Value = 3.24
$ws.Cells.Item(8, 15).Value = '22/07/2023 14:42'
$ws.Cells.Item(8, 16).Value = 3.41
$ws.Cells.Item(8, 17).Value = '22/07/2023 17:59'
$ws.Cells.Item(8, 18).Value = 2.57
$ws.Cells.Item(8, 19).Value = '22/07/2023 14:42'
$ws.Cells.Item(8, 20).Value = 2.42
$ws.Cells.Item(8, 21).Value = '22/07/2023 17:59'
$ws.Cells.Item(8, 22).Value = 'https://www.betexplorer.com/football/poland/division-2/gks-jastrzebie-stal-stalowa-wola/84hLxcUH/'
# --- Rows 26,27: swap ---
$ws.Cells.Item(26, 6).Value = 'Hutnik Krakow'
$ws.Cells.Item(26, 7).Value = 1
$ws.Cells.Item(26, 8).Value = 'Zaglebie II'
$ws.Cells.Item(26, 9).Value = 1
$ws.Cells.Item(26, 10).Value = 1.89
$ws.Cells.Item(26, 11).Value = '05/08/2023 22:29'
$ws.Cells.Item(26, 12).Value = 1.94
$ws.Cells.Item(26, 13).Value = '06/08/2023 16:50'
$ws.Cells.Item(26, 14).Value = 3.66
$ws.Cells.Item(26, 15).Value = '05/08/2023 22:29'
$ws.Cells.Item(26, 16).Value = 3.54
$ws.Cells.Item(26, 17).Value = '06/08/2023 16:50'
$ws.Cells.Item(26, 18).Value = 3.62
$ws.Cells.Item(26, 19).Value = '05/08/2023 22:29'
$ws.Cells.Item(26, 20).Value = 3.65
$ws.Cells.Item(26, 21).Value = '06/08/2023 16:50'
$ws.Cells.Item(26, 22).Value = 'https://www.betexplorer.com/football/poland/division-2/hutnik-krakow-zaglebie/xdUUh9KF/'
$ws.Cells.Item(27, 6).Value = 'Stezyca'
$ws.Cells.Item(27, 7).Value = 1
$ws.Cells.Item(27, 8).Value = 'Sandecja Nowy S.'
$ws.Cells.Item(27, 9).Value = 1
$ws.Cells.Item(27, 10).Value = 2.27
$ws.Cells.Item(27, 11).Value = '05/08/2023 22:29'
$ws.Cells.Item(27, 12).Value = 2.32
$ws.Cells.Item(27, 13).Value = '06/08/2023 16:44'
$ws.Cells.Item(27, 14).Value = 3.23
$ws.Cells.Item(27, 15).Value = '05/08/2023 22:29'
$ws.Cells.Item(27, 16).Value = 3.31
$ws.Cells.Item(27, 17).Value = '06/08/2023 15:54'
$ws.Cells.Item(27, 18).Value = 3.05
$ws.Cells.Item(27, 19).Value = '05/08/2023 22:29'
$ws.Cells.Item(27, 20).Value = 2.96
$ws.Cells.Item(27, 21).Value = '06/08/2023 16:44'
$ws.Cells.Item(27, 22).Value = 'https://www.betexplorer.com/football/poland/division-2/stezyca-sandecja-nowy-s/WITYiTZL/'
# --- Rows 41,42: swap ---
$ws.Cells.Item(41, 6).Value = 'Olimpia Elblag'
$ws.Cells.Item(41, 7).Value = 0
$ws.Cells.Item(41, 8).Value = 'Stezyca'
$ws.Cells.Item(41, 9).Value = 1
$ws.Cells.Item(41, 10).Value = 2.27
$ws.Cells.Item(41, 11).Value = '19/08/2023 08:43'
$ws.Cells.Item(41, 12).Value = 2.27
$ws.Cells.Item(41, 13).Value = '19/08/2023 08:43'
$ws.Cells.Item(41, 14).Value = 3.18
$ws.Cells.Item(41, 15).Value = '19/08/2023 08:43'
$ws.Cells.Item(41, 16).Value = 3.2
$ws.Cells.Item(41, 17).Value = '19/08/2023 16:05'
$ws.Cells.Item(41, 18).Value = 3.09
$ws.Cells.Item(41, 19).Value = '19/08/2023 08:43'
$ws.Cells.Item(41, 20).Value = 3.09
$ws.Cells.Item(41, 21).Value = '19/08/2023 08:43'
$ws.Cells.Item(41, 22).Value = 'https://www.betexplorer.com/football/poland/division-2/olimpia-elblag-stezyca/OYTnIPmR/'
$ws.Cells.Item(42, 6).Value = 'Kotwica Kolobrzeg'
$ws.Cells.Item(42, 7).Value = 1
$ws.Cells.Item(42, 8).Value = 'Ol. Grudziadz'
$ws.Cells.Item(42, 9).Value = 2
$ws.Cells.Item(42, 10).Value = 1.91
$ws.Cells.Item(42, 11).Value = '19/08/2023 08:43'
$ws.Cells.Item(42, 12).Value = 2.08
$ws.Cells.Item(42, 13).Value = '19/08/2023 17:46'
$ws.Cells.Item(42, 14).Value = 3.45
$ws.Cells.Item(42, 15).Value = '19/08/2023 08:43'
$ws.Cells.Item(42, 16).Value = 3.27
$ws.Cells.Item(42, 17).Value = '19/08/2023 17:46'
$ws.Cells.Item(42, 18).Value = 3.77
$ws.Cells.Item(42, 19).Value = '19/08/2023 08:43'
$ws.Cells.Item(42, 20).Value = 3.51
$ws.Cells.Item(42, 21).Value = '19/08/2023 17:46'
$ws.Cells.Item(42, 22).Value = 'https://www.betexplorer.com/football/poland/division-2/kotwica-kolobrzeg-ol-grudziadz/KIUHC3Bl/'
# --- Rows 51,52: swap ---
$ws.Cells.Item(51, 6).Value = 'Pogon Siedlce'
$ws.Cells.Item(51, 7).Value = 3
$ws.Cells.Item(51, 8).Value = 'Sandecja Nowy S.'
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 10).Value = 2.28
$ws.Cells.Item(51, 11).Value = '26/08/2023 13:13'
$ws.Cells.Item(51, 12).Value = 2.77
$ws.Cells.Item(51, 13).Value = '26/08/2023 16:56'
$ws.Cells.Item(51, 14).Value = 3.08
$ws.Cells.Item(51, 15).Value = '26/08/2023 13:13'
$ws.Cells.Item(51, 16).Value = 3.12
$ws.Cells.Item(51, 17).Value = '26/08/2023 16:55'
$ws.Cells.Item(51, 18).Value = 3.09
$ws.Cells.Item(51, 19).Value = '26/08/2023 13:13'
$ws.Cells.Item(51, 20).Value = 2.56
$ws.Cells.Item(51, 21).Value = '26/08/2023 16:56'
$ws.Cells.Item(51, 22).Value = 'https://www.betexplorer.com/football/poland/division-2/pogon-siedlce-sandecja-nowy-s/bT3PVrIQ/'
$ws.Cells.Item(52, 6).Value = 'GKS Jastrzebie'
$ws.Cells.Item(52, 7).Value = 4
$ws.Cells.Item(52, 8).Value = 'Polonia Bytom'
$ws.Cells.Item(52, 9).Value = 2
$ws.Cells.Item(52, 10).Value = 1.85
$ws.Cells.Item(52, 11).Value = '26/08/2023 13:13'
$ws.Cells.Item(52, 12).Value = 1.83
$ws.Cells.Item(52, 13).Value = '26/08/2023 16:58'
$ws.Cells.Item(52, 14).Value = 3.42
$ws.Cells.Item(52, 15).Value = '26/08/2023 13:13'
$ws.Cells.Item(52, 16).Value = 3.6
$ws.Cells.Item(52, 17).Value = '26/08/2023 16:58'
$ws.Cells.Item(52, 18).Value = 4.09
$ws.Cells.Item(52, 19).Value = '26/08/2023 13:13'
$ws.Cells.Item(52, 20).Value = 4.03
$ws.Cells.Item(52, 21).Value = '26/08/2023 16:58'
$ws.Cells.Item(52, 22).Value = 'https://www.betexplorer.com/football/poland/division-2/gks-jastrzebie-polonia-bytom/464LWO2K/'
# --- Rows 91,93: swap ---
$ws.Cells.Item(91, 6).Value = 'Sandecja Nowy S.'
$ws.Cells.Item(91, 7).Value = 0
$ws.Cells.Item(91, 8).Value = 'KKS Kalisz'
$ws.Cells.Item(91, 9).Value = 2
$ws.Cells.Item(91, 10).Value = 2.75
$ws.Cells.Item(91, 11).Value = '29/09/2023 02:12'
$ws.Cells.Item(91, 12).Value = 2.78
$ws.Cells.Item(91, 13).Value = '30/09/2023 14:41'
$ws.Cells.Item(91, 14).Value = 3.09
$ws.Cells.Item(91, 15).Value = '29/09/2023 02:12'
$ws.Cells.Item(91, 16).Value = 3.12
$ws.Cells.Item(91, 17).Value = '30/09/2023 14:41'
$ws.Cells.Item(91, 18).Value = 2.35
$ws.Cells.Item(91, 19).Value = '29/09/2023 02:12'
$ws.Cells.Item(91, 20).Value = 2.56
$ws.Cells.Item(91, 21).Value = '30/09/2023 14:41'
$ws.Cells.Item(91, 22).Value = 'https://www.betexplorer.com/football/poland/division-2/sandecja-nowy-s-kks-kalisz/KbAXndAF/'
$ws.Cells.Item(93, 6).Value = 'Hutnik Krakow'
$ws.Cells.Item(93, 7).Value = 0
$ws.Cells.Item(93, 8).Value = 'Chojniczanka'
$ws.Cells.Item(93, 9).Value = 2
$ws.Cells.Item(93, 10).Value = 2.06
$ws.Cells.Item(93, 11).Value = '29/09/2023 02:12'
$ws.Cells.Item(93, 12).Value = 2.24
$ws.Cells.Item(93, 13).Value = '30/09/2023 14:43'
$ws.Cells.Item(93, 14).Value = 3.21
$ws.Cells.Item(93, 15).Value = '29/09/2023 02:12'
$ws.Cells.Item(93, 16).Value = 3.49
$ws.Cells.Item(93, 17).Value = '30/09/2023 14:41'
$ws.Cells.Item(93, 18).Value = 3.14
$ws.Cells.Item(93, 19).Value = '29/09/2023 02:12'
$ws.Cells.Item(93, 20).Value = 2.95
$ws.Cells.Item(93, 21).Value = '30/09/2023 14:43'
$ws.Cells.Item(93, 22).Value = 'https://www.betexplorer.com/football/poland/division-2/hutnik-krakow-chojniczanka/t8BTmxe9/'
# --- Rows 99,100: swap ---
$ws.Cells.Item(99, 6).Value = 'LKS Lodz II'
$ws.Cells.Item(99, 7).Value = 2
$ws.Cells.Item(99, 8).Value = 'Kotwica Kolobrzeg'
$ws.Cells.Item(99, 9).Value = 3
$ws.Cells.Item(99, 10).Value = 2.3
$ws.Cells.Item(99, 11).Value = '05/10/2023 05:12'
$ws.Cells.Item(99, 12).Value = 2.61
$ws.Cells.Item(99, 13).Value = '06/10/2023 17:55'
$ws.Cells.Item(99, 14).Value = 3.18
$ws.Cells.Item(99, 15).Value = '05/10/2023 05:12'
$ws.Cells.Item(99, 16).Value = 3.67
$ws.Cells.Item(99, 17).Value = '06/10/2023 17:55'
$ws.Cells.Item(99, 18).Value = 2.74
$ws.Cells.Item(99, 19).Value = '05/10/2023 05:12'
$ws.Cells.Item(99, 20).Value = 2.41
$ws.Cells.Item(99, 21).Value = '06/10/2023 17:52'
$ws.Cells.Item(99, 22).Value = 'https://www.betexplorer.com/football/poland/division-2/lks-lodz-kotwica-kolobrzeg/fgIqtog0/'
$ws.Cells.Item(100, 6).Value = 'S. Wola'
$ws.Cells.Item(100, 7).Value = 3
$ws.Cells.Item(100, 8).Value = 'Lech Poznan II'
$ws.Cells.Item(100, 9).Value = 1
$ws.Cells.Item(100, 10).Value = 1.7
$ws.Cells.Item(100, 11).Value = '05/10/2023 05:12'
$ws.Cells.Item(100, 12).Value = 1.32
$ws.Cells.Item(100, 13).Value = '06/10/2023 17:49'
$ws.Cells.Item(100, 14).Value = 3.61
$ws.Cells.Item(100, 15).Value = '05/10/2023 05:12'
$ws.Cells.Item(100, 16).Value = 5.14
$ws.Cells.Item(100, 17).Value = '06/10/2023 17:57'
$ws.Cells.Item(100, 18).Value = 4
$ws.Cells.Item(100, 19).Value = '05/10/2023 05:12'
$ws.Cells.Item(100, 20).Value = 8.6
$ws.Cells.Item(100, 21).Value = '06/10/2023 17:57'
$ws.Cells.Item(100, 22).Value = 'https://www.betexplorer.com/football/poland/division-2/stal-stalowa-wola-lech-poznan/6TLmu586/'
# --- Rows 148,149: swap ---
$ws.Cells.Item(148, 6).Value = 'Zaglebie II'
$ws.Cells.Item(148, 7).Value = 4
$ws.Cells.Item(148, 8).Value = 'S. Wola'
$ws.Cells.Item(148, 9).Value = 0
$ws.Cells.Item(148, 10).Value = 2.3
$ws.Cells.Item(148, 11).Value = '11/11/2023 01:13'
$ws.Cells.Item(148, 12).Value = 2.72
$ws.Cells.Item(148, 13).Value = '12/11/2023 12:51'
$ws.Cells.Item(148, 14).Value = 3.19
$ws.Cells.Item(148, 15).Value = '11/11/2023 01:13'
$ws.Cells.Item(148, 16).Value = 3.15
$ws.Cells.Item(148, 17).Value = '12/11/2023 12:51'
$ws.Cells.Item(148, 18).Value = 2.81
$ws.Cells.Item(148, 19).Value = '11/11/2023 01:13'
$ws.Cells.Item(148, 20).Value = 2.59
$ws.Cells.Item(148, 21).Value = '12/11/2023 12:51'
$ws.Cells.Item(148, 22).Value = 'https://www.betexplorer.com/football/poland/division-2/zaglebie-stal-stalowa-wola/zqWeXABL/'
$ws.Cells.Item(149, 6).Value = 'Stezyca'
$ws.Cells.Item(149, 7).Value = 1
$ws.Cells.Item(149, 8).Value = 'Stomil Olsztyn'
$ws.Cells.Item(149, 9).Value = 0
$ws.Cells.Item(149, 10).Value = 1.95
$ws.Cells.Item(149, 11).Value = '11/11/2023 01:13'
$ws.Cells.Item(149, 12).Value = 1.9
$ws.Cells.Item(149, 13).Value = '12/11/2023 12:39'
$ws.Cells.Item(149, 14).Value = 3.28
$ws.Cells.Item(149, 15).Value = '11/11/2023 01:13'
$ws.Cells.Item(149, 16).Value = 3.33
$ws.Cells.Item(149, 17).Value = '12/11/2023 12:39'
$ws.Cells.Item(149, 18).Value = 3.35
$ws.Cells.Item(149, 19).Value = '11/11/2023 01:13'
$ws.Cells.Item(149, 20).Value = 4.08
$ws.Cells.Item(149, 21).Value = '12/11/2023 12:39'
$ws.Cells.Item(149, 22).Value = 'https://www.betexplorer.com/football/poland/division-2/stezyca-stomil-olsztyn/dhXiYjdF/'
# --- Rows 153,154: swap ---
$ws.Cells.Item(153, 6).Value = 'Sandecja Nowy S.'
$ws.Cells.Item(153, 7).Value = 1
$ws.Cells.Item(153, 8).Value = 'Ol. Grudziadz'
$ws.Cells.Item(153, 9).Value = 0
$ws.Cells.Item(153, 10).Value = 2.36
$ws.Cells.Item(153, 11).Value = '17/11/2023 01:12'
$ws.Cells.Item(153, 12).Value = 2.44
$ws.Cells.Item(153, 13).Value = '18/11/2023 11:11'
$ws.Cells.Item(153, 14).Value = 3.14
$ws.Cells.Item(153, 15).Value = '17/11/2023 01:12'
$ws.Cells.Item(153, 16).Value = 3.07
$ws.Cells.Item(153, 17).Value = '18/11/2023 11:11'
$ws.Cells.Item(153, 18).Value = 2.69
$ws.Cells.Item(153, 19).Value = '17/11/2023 01:12'
$ws.Cells.Item(153, 20).Value = 2.98
$ws.Cells.Item(153, 21).Value = '18/11/2023 11:11'
$ws.Cells.Item(153, 22).Value = 'https://www.betexplorer.com/football/poland/division-2/sandecja-nowy-s-ol-grudziadz/WAozoUuL/'
$ws.Cells.Item(154, 6).Value = 'Zaglebie II'
$ws.Cells.Item(154, 7).Value = 2
$ws.Cells.Item(154, 8).Value = 'Stomil Olsztyn'
$ws.Cells.Item(154, 9).Value = 1
$ws.Cells.Item(154, 10).Value = 2.19
$ws.Cells.Item(154, 11).Value = '17/11/2023 01:12'
$ws.Cells.Item(154, 12).Value = 2.13
$ws.Cells.Item(154, 13).Value = '18/11/2023 12:51'
$ws.Cells.Item(154, 14).Value = 3.31
$ws.Cells.Item(154, 15).Value = '17/11/2023 01:12'
$ws.Cells.Item(154, 16).Value = 3.48
$ws.Cells.Item(154, 17).Value = '18/11/2023 12:51'
$ws.Cells.Item(154, 18).Value = 2.8
$ws.Cells.Item(154, 19).Value = '17/11/2023 01:12'
$ws.Cells.Item(154, 20).Value = 3.18
$ws.Cells.Item(154, 21).Value = '18/11/2023 12:51'
$ws.Cells.Item(154, 22).Value = 'https://www.betexplorer.com/football/poland/division-2/zaglebie-stomil-olsztyn/tpyAvSIl/'
# --- New rows 161-167 (match indices 160-166) ---
# Row 161: Indice=160
$ws.Cells.Item(161, 1).Value = 160
$ws.Cells.Item(160, 1).Copy()
$ws.Cells.Item(161, 1).PasteSpecial(-4122)
$ws.Cells.Item(161, 2).Value = 'poland'
$ws.Cells.Item(161, 3).Value = 'division-2'
$ws.Cells.Item(161, 4).Value = '2023-2024'
$ws.Cells.Item(161, 5).Value = 45255.5
$ws.Cells.Item(160, 5).Copy()
$ws.Cells.Item(161, 5).PasteSpecial(-4122)
$ws.Cells.Item(161, 6).Value = 'Lech Poznan II'
$ws.Cells.Item(161, 7).Value = 1
$ws.Cells.Item(161, 8).Value = 'Chojniczanka'
$ws.Cells.Item(161, 9).Value = 2
$ws.Cells.Item(161, 10).Value = 3.02
$ws.Cells.Item(161, 11).Value = '24/11/2023 00:13'
$ws.Cells.Item(161, 12).Value = 3.53
$ws.Cells.Item(161, 13).Value = '25/11/2023 11:59'
$ws.Cells.Item(161, 14).Value = 3.37
$ws.Cells.Item(161, 15).Value = '24/11/2023 00:13'
$ws.Cells.Item(161, 16).Value = 3.51
$ws.Cells.Item(161, 17).Value = '25/11/2023 11:59'
$ws.Cells.Item(161, 18).Value = 2.1
$ws.Cells.Item(161, 19).Value = '24/11/2023 00:13'
$ws.Cells.Item(161, 20).Value = 1.99
$ws.Cells.Item(161, 21).Value = '25/11/2023 11:59'
$ws.Cells.Item(161, 22).Value = 'https://www.betexplorer.com/football/poland/division-2/lech-poznan-chojniczanka/APSQzpID/'
$ws.Cells.Item(161, 5).NumberFormat = 'YYYY-MM-DD HH:MM:SS'

# Row 162: Indice=161
$ws.Cells.Item(162, 1).Value = 161
$ws.Cells.Item(160, 1).Copy()
$ws.Cells.Item(162, 1).PasteSpecial(-4122)
$ws.Cells.Item(162, 2).Value = 'poland'
$ws.Cells.Item(162, 3).Value = 'division-2'
$ws.Cells.Item(162, 4).Value = '2023-2024'
$ws.Cells.Item(162, 5).Value = 45255.52083333334
$ws.Cells.Item(160, 5).Copy()
$ws.Cells.Item(162, 5).PasteSpecial(-4122)
$ws.Cells.Item(162, 6).Value = 'Hutnik Krakow'
$ws.Cells.Item(162, 7).Value = 3
$ws.Cells.Item(162, 8).Value = 'Stezyca'
$ws.Cells.Item(162, 9).Value = 1
$ws.Cells.Item(162, 10).Value = 2.04
$ws.Cells.Item(162, 11).Value = '24/11/2023 00:42'
$ws.Cells.Item(162, 12).Value = 2.15
$ws.Cells.Item(162, 13).Value = '25/11/2023 12:11'
$ws.Cells.Item(162, 14).Value = 3.26
$ws.Cells.Item(162, 15).Value = '24/11/2023 00:42'
$ws.Cells.Item(162, 16).Value = 3.33
$ws.Cells.Item(162, 17).Value = '25/11/2023 12:11'
$ws.Cells.Item(162, 18).Value = 3.15
$ws.Cells.Item(162, 19).Value = '24/11/2023 00:42'
$ws.Cells.Item(162, 20).Value = 3.28
$ws.Cells.Item(162, 21).Value = '25/11/2023 12:11'
$ws.Cells.Item(162, 22).Value = 'https://www.betexplorer.com/football/poland/division-2/hutnik-krakow-stezyca/MJAZxzas/'
$ws.Cells.Item(162, 5).NumberFormat = 'YYYY-MM-DD HH:MM:SS'

# Row 163: Indice=162
$ws.Cells.Item(163, 1).Value = 162
$ws.Cells.Item(160, 1).Copy()
$ws.Cells.Item(163, 1).PasteSpecial(-4122)
$ws.Cells.Item(163, 2).Value = 'poland'
$ws.Cells.Item(163, 3).Value = 'division-2'
$ws.Cells.Item(163, 4).Value = '2023-2024'
$ws.Cells.Item(163, 5).Value = 45255.54166666666
$ws.Cells.Item(160, 5).Copy()
$ws.Cells.Item(163, 5).PasteSpecial(-4122)
$ws.Cells.Item(163, 6).Value = 'Sandecja Nowy S.'
$ws.Cells.Item(163, 7).Value = 2
$ws.Cells.Item(163, 8).Value = 'Skra'
$ws.Cells.Item(163, 9).Value = 2
$ws.Cells.Item(163, 10).Value = 2.46
$ws.Cells.Item(163, 11).Value = '24/11/2023 01:13'
$ws.Cells.Item(163, 12).Value = 2.94
$ws.Cells.Item(163, 13).Value = '25/11/2023 12:44'
$ws.Cells.Item(163, 14).Value = 3.02
$ws.Cells.Item(163, 15).Value = '24/11/2023 01:13'
$ws.Cells.Item(163, 16).Value = 2.91
$ws.Cells.Item(163, 17).Value = '25/11/2023 12:41'
$ws.Cells.Item(163, 18).Value = 2.66
$ws.Cells.Item(163, 19).Value = '24/11/2023 01:13'
$ws.Cells.Item(163, 20).Value = 2.58
$ws.Cells.Item(163, 21).Value = '25/11/2023 12:44'
$ws.Cells.Item(163, 22).Value = 'https://www.betexplorer.com/football/poland/division-2/sandecja-nowy-s-skra-czestochowa/pbJ7sGjQ/'
$ws.Cells.Item(163, 5).NumberFormat = 'YYYY-MM-DD HH:MM:SS'

# Row 164: Indice=163
$ws.Cells.Item(164, 1).Value = 163
$ws.Cells.Item(160, 1).Copy()
$ws.Cells.Item(164, 1).PasteSpecial(-4122)
$ws.Cells.Item(164, 2).Value = 'poland'
$ws.Cells.Item(164, 3).Value = 'division-2'
$ws.Cells.Item(164, 4).Value = '2023-2024'
$ws.Cells.Item(164, 5).Value = 45255.54166666666
$ws.Cells.Item(160, 5).Copy()
$ws.Cells.Item(164, 5).PasteSpecial(-4122)
$ws.Cells.Item(164, 6).Value = 'Stomil Olsztyn'
$ws.Cells.Item(164, 7).Value = 2
$ws.Cells.Item(164, 8).Value = 'GKS Jastrzebie'
$ws.Cells.Item(164, 9).Value = 0
$ws.Cells.Item(164, 10).Value = 2.32
$ws.Cells.Item(164, 11).Value = '24/11/2023 01:13'
$ws.Cells.Item(164, 12).Value = 2.15
$ws.Cells.Item(164, 13).Value = '25/11/2023 12:57'
$ws.Cells.Item(164, 14).Value = 3.17
$ws.Cells.Item(164, 15).Value = '24/11/2023 01:13'
$ws.Cells.Item(164, 16).Value = 3.48
$ws.Cells.Item(164, 17).Value = '25/11/2023 12:57'
$ws.Cells.Item(164, 18).Value = 2.8
$ws.Cells.Item(164, 19).Value = '24/11/2023 01:13'
$ws.Cells.Item(164, 20).Value = 3.14
$ws.Cells.Item(164, 21).Value = '25/11/2023 12:57'
$ws.Cells.Item(164, 22).Value = 'https://www.betexplorer.com/football/poland/division-2/stomil-olsztyn-gks-jastrzebie/rZbryETg/'
$ws.Cells.Item(164, 5).NumberFormat = 'YYYY-MM-DD HH:MM:SS'

# Row 165: Indice=164
$ws.Cells.Item(165, 1).Value = 164
$ws.Cells.Item(160, 1).Copy()
$ws.Cells.Item(165, 1).PasteSpecial(-4122)
$ws.Cells.Item(165, 2).Value = 'poland'
$ws.Cells.Item(165, 3).Value = 'division-2'
$ws.Cells.Item(165, 4).Value = '2023-2024'
$ws.Cells.Item(165, 5).Value = 45255.66666666666
$ws.Cells.Item(160, 5).Copy()
$ws.Cells.Item(165, 5).PasteSpecial(-4122)
$ws.Cells.Item(165, 6).Value = 'KKS Kalisz'
$ws.Cells.Item(165, 7).Value = 3
$ws.Cells.Item(165, 8).Value = 'Zaglebie II'
$ws.Cells.Item(165, 9).Value = 1
$ws.Cells.Item(165, 10).Value = 1.68
$ws.Cells.Item(165, 11).Value = '24/11/2023 04:12'
$ws.Cells.Item(165, 12).Value = 1.81
$ws.Cells.Item(165, 13).Value = '25/11/2023 15:51'
$ws.Cells.Item(165, 14).Value = 3.79
$ws.Cells.Item(165, 15).Value = '24/11/2023 04:12'
$ws.Cells.Item(165, 16).Value = 3.93
$ws.Cells.Item(165, 17).Value = '25/11/2023 15:51'
$ws.Cells.Item(165, 18).Value = 3.86
$ws.Cells.Item(165, 19).Value = '24/11/2023 04:12'
$ws.Cells.Item(165, 20).Value = 3.79
$ws.Cells.Item(165, 21).Value = '25/11/2023 15:51'
$ws.Cells.Item(165, 22).Value = 'https://www.betexplorer.com/football/poland/division-2/kks-kalisz-zaglebie/4r0vxfEm/'
$ws.Cells.Item(165, 5).NumberFormat = 'YYYY-MM-DD HH:MM:SS'

# Row 166: Indice=165
$ws.Cells.Item(166, 1).Value = 165
$ws.Cells.Item(160, 1).Copy()
$ws.Cells.Item(166, 1).PasteSpecial(-4122)
$ws.Cells.Item(166, 2).Value = 'poland'
$ws.Cells.Item(166, 3).Value = 'division-2'
$ws.Cells.Item(166, 4).Value = '2023-2024'
$ws.Cells.Item(166, 5).Value = 45255.70833333334
$ws.Cells.Item(160, 5).Copy()
$ws.Cells.Item(166, 5).PasteSpecial(-4122)
$ws.Cells.Item(166, 6).Value = 'Polonia Bytom'
$ws.Cells.Item(166, 7).Value = 2
$ws.Cells.Item(166, 8).Value = 'Ol. Grudziadz'
$ws.Cells.Item(166, 9).Value = 1
$ws.Cells.Item(166, 10).Value = 2.26
$ws.Cells.Item(166, 11).Value = '24/11/2023 05:13'
$ws.Cells.Item(166, 12).Value = 2.32
$ws.Cells.Item(166, 13).Value = '25/11/2023 16:59'
$ws.Cells.Item(166, 14).Value = 3.18
$ws.Cells.Item(166, 15).Value = '24/11/2023 05:13'
$ws.Cells.Item(166, 16).Value = 3.24
$ws.Cells.Item(166, 17).Value = '25/11/2023 16:59'
$ws.Cells.Item(166, 18).Value = 2.8
$ws.Cells.Item(166, 19).Value = '24/11/2023 05:13'
$ws.Cells.Item(166, 20).Value = 3.02
$ws.Cells.Item(166, 21).Value = '25/11/2023 15:13'
$ws.Cells.Item(166, 22).Value = 'https://www.betexplorer.com/football/poland/division-2/polonia-bytom-ol-grudziadz/jyK3rdyK/'
$ws.Cells.Item(166, 5).NumberFormat = 'YYYY-MM-DD HH:MM:SS'

# Row 167: Indice=166
$ws.Cells.Item(167, 1).Value = 166
$ws.Cells.Item(160, 1).Copy()
$ws.Cells.Item(167, 1).PasteSpecial(-4122)
$ws.Cells.Item(167, 2).Value = 'poland'
$ws.Cells.Item(167, 3).Value = 'division-2'
$ws.Cells.Item(167, 4).Value = '2023-2024'
$ws.Cells.Item(167, 5).Value = 45256.72916666666
$ws.Cells.Item(160, 5).Copy()
$ws.Cells.Item(167, 5).PasteSpecial(-4122)
$ws.Cells.Item(167, 6).Value = 'LKS Lodz II'
$ws.Cells.Item(167, 7).Value = 2
$ws.Cells.Item(167, 8).Value = 'Wisla Pulawy'
$ws.Cells.Item(167, 9).Value = 1
$ws.Cells.Item(167, 10).Value = 2.16
$ws.Cells.Item(167, 11).Value = '25/11/2023 05:43'
$ws.Cells.Item(167, 12).Value = 2.18
$ws.Cells.Item(167, 13).Value = '26/11/2023 17:26'
$ws.Cells.Item(167, 14).Value = 3.42
$ws.Cells.Item(167, 15).Value = '25/11/2023 05:43'
$ws.Cells.Item(167, 16).Value = 3.59
$ws.Cells.Item(167, 17).Value = '26/11/2023 17:26'
$ws.Cells.Item(167, 18).Value = 2.79
$ws.Cells.Item(167, 19).Value = '25/11/2023 05:43'
$ws.Cells.Item(167, 20).Value = 2.98
$ws.Cells.Item(167, 21).Value = '26/11/2023 17:26'
$ws.Cells.Item(167, 22).Value = 'https://www.betexplorer.com/football/poland/division-2/lks-lodz-wisla-pulawy/ja2jZib6/'
$ws.Cells.Item(167, 5).NumberFormat = 'YYYY-MM-DD HH:MM:SS'

